$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("general")

$ws.Range("B35").Value = "Can you list the organizations that are available in PLM?"
$ws.Range("B36").Value = "What kinds of organizations can students join in PLM?"
$ws.Range("B37").Value = "Are there any specific organizations that are popular among PLM students?"
$ws.Range("B38").Value = "Could you provide some information on the different organizations that exist in PLM?"
$ws.Range("B40").Value = "Can you tell me about the various organizations that operate in PLM?"
$ws.Range("B41").Value = "Which organizations are currently active in PLM?"
$ws.Range("B42").Value = "How many organizations are there in PLM, and what are they?"
